$wb = $excel.ActiveWorkbook

# Belgium is the template for the new "Czech" sheet (same layout/styles).
$belgium = $wb.Worksheets.Item("Belgium")

# Copy Belgium and place the copy right after it -> becomes the 3rd sheet.
$belgium.Copy($null, $belgium)
$czech = $wb.Worksheets.Item(3)
$czech.Name = "Czech"

# Fill in the Czech market test data.
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3478/T1732"

# Columns re-sized (narrower text than Belgium's) to match the authored sheet.
$czech.Columns.Item(2).ColumnWidth = 35.42
$czech.Columns.Item(3).ColumnWidth = 12.92
$czech.Columns.Item(4).ColumnWidth = 17.25

# Belgium's selection settles back to the full data range (no longer the active tab).
$belgium.Range("A1:D11").Select()

# Czech becomes the active sheet/tab, with C4 selected.
$czech.Activate()
$czech.Range("C4").Select()
